$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.272596120834351
$ws.Range("B1").Value = 2.488819122314453
$ws.Range("C1").Value = 3.516110897064209
$ws.Range("D1").Value = 3.043841600418091
$ws.Range("E1").Value = 1.068212032318115
